$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Admin2"
$ws.Range("B3").Value = "admin456"
$ws.Range("A4").Value = "Admin3"
$ws.Range("B4").Value = "admin789"

$ws.Range("B7").Select()
